# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# Row 21: F21  1337/1338 -> 1362
# Row 25: F25  1932      -> 1934
# Row 29: F29  191       -> 192

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F21").Value = 1362
    $ws.Range("F25").Value = 1934
    $ws.Range("F29").Value = 192
}
